$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing rows 7-83: shift weekly data down by one row, with a new week inserted at row 7
$ws.Range("D7").Value = 44552
$ws.Range("J7").Value = 106
$ws.Range("K7").Value = 7000
$ws.Range("L7").Value = 8000
$ws.Range("M7").Value = 7500
$ws.Range("P7").Value = 375

$ws.Range("D8").Value = 44545
$ws.Range("J8").Value = 160
$ws.Range("K8").Value = 6000
$ws.Range("L8").Value = 7000
$ws.Range("M8").Value = 6500
$ws.Range("P8").Value = 325

$ws.Range("D9").Value = 44421
$ws.Range("J9").Value = 180
$ws.Range("K9").Value = 7000
$ws.Range("L9").Value = 8000
$ws.Range("M9").Value = 7500
$ws.Range("P9").Value = 375

$ws.Range("D10").Value = 44314
$ws.Range("J10").Value = 160
$ws.Range("K10").Value = 8000
$ws.Range("L10").Value = 8000
$ws.Range("M10").Value = 8000
$ws.Range("P10").Value = 400

$ws.Range("D11").Value = 44229
$ws.Range("J11").Value = 50
$ws.Range("K11").Value = 7000
$ws.Range("L11").Value = 7000
$ws.Range("M11").Value = 7000
$ws.Range("P11").Value = 350

$ws.Range("D12").Value = 44503
$ws.Range("J12").Value = 97
$ws.Range("K12").Value = 7000
$ws.Range("L12").Value = 8000
$ws.Range("M12").Value = 7505
$ws.Range("P12").Value = 375

$ws.Range("D13").Value = 44497
$ws.Range("J13").Value = 180
$ws.Range("K13").Value = 6000
$ws.Range("L13").Value = 7000
$ws.Range("M13").Value = 6556
$ws.Range("P13").Value = 328

$ws.Range("D14").Value = 44435
$ws.Range("J14").Value = 302
$ws.Range("K14").Value = 7000
$ws.Range("L14").Value = 8000
$ws.Range("M14").Value = 7500
$ws.Range("P14").Value = 375

$ws.Range("D15").Value = 44295
$ws.Range("J15").Value = 70
$ws.Range("K15").Value = 8000
$ws.Range("L15").Value = 8000
$ws.Range("M15").Value = 8000
$ws.Range("P15").Value = 400

$ws.Range("D16").Value = 44356
$ws.Range("J16").Value = 160
$ws.Range("K16").Value = 7000
$ws.Range("L16").Value = 8000
$ws.Range("M16").Value = 7500
$ws.Range("P16").Value = 375

$ws.Range("D17").Value = 44244
$ws.Range("J17").Value = 70
$ws.Range("K17").Value = 8000
$ws.Range("L17").Value = 8000
$ws.Range("M17").Value = 8000
$ws.Range("P17").Value = 400

$ws.Range("D18").Value = 44342
$ws.Range("J18").Value = 160
$ws.Range("K18").Value = 8000
$ws.Range("L18").Value = 8000
$ws.Range("M18").Value = 8000
$ws.Range("P18").Value = 400

$ws.Range("D19").Value = 44524
$ws.Range("J19").Value = 160
$ws.Range("K19").Value = 6000
$ws.Range("L19").Value = 7000
$ws.Range("M19").Value = 6500
$ws.Range("P19").Value = 325

$ws.Range("D20").Value = 44189
$ws.Range("J20").Value = 50
$ws.Range("K20").Value = 8000
$ws.Range("L20").Value = 8000
$ws.Range("M20").Value = 8000
$ws.Range("P20").Value = 400

$ws.Range("D21").Value = 44442
$ws.Range("J21").Value = 180
$ws.Range("K21").Value = 7000
$ws.Range("L21").Value = 8000
$ws.Range("M21").Value = 7500
$ws.Range("P21").Value = 375

$ws.Range("D22").Value = 44363
$ws.Range("J22").Value = 160
$ws.Range("K22").Value = 8000
$ws.Range("L22").Value = 8000
$ws.Range("M22").Value = 8000
$ws.Range("P22").Value = 400

$ws.Range("D23").Value = 44215
$ws.Range("J23").Value = 80
$ws.Range("K23").Value = 7000
$ws.Range("L23").Value = 7000
$ws.Range("M23").Value = 7000
$ws.Range("P23").Value = 350

$ws.Range("D24").Value = 44517
$ws.Range("J24").Value = 160
$ws.Range("K24").Value = 6000
$ws.Range("L24").Value = 7000
$ws.Range("M24").Value = 6500
$ws.Range("P24").Value = 325

$ws.Range("D25").Value = 44484
$ws.Range("J25").Value = 160
$ws.Range("K25").Value = 7000
$ws.Range("L25").Value = 8000
$ws.Range("M25").Value = 7500
$ws.Range("P25").Value = 375

$ws.Range("D26").Value = 44426
$ws.Range("J26").Value = 97
$ws.Range("K26").Value = 7000
$ws.Range("L26").Value = 8000
$ws.Range("M26").Value = 7505
$ws.Range("P26").Value = 375

$ws.Range("D27").Value = 44358
$ws.Range("J27").Value = 160
$ws.Range("K27").Value = 7500
$ws.Range("L27").Value = 8000
$ws.Range("M27").Value = 7750
$ws.Range("P27").Value = 388

$ws.Range("D28").Value = 44477
$ws.Range("J28").Value = 160
$ws.Range("K28").Value = 7000
$ws.Range("L28").Value = 8000
$ws.Range("M28").Value = 7500
$ws.Range("P28").Value = 375

$ws.Range("D29").Value = 44203
$ws.Range("J29").Value = 50
$ws.Range("K29").Value = 7000
$ws.Range("L29").Value = 8000
$ws.Range("M29").Value = 7400
$ws.Range("P29").Value = 370

$ws.Range("D30").Value = 44384
$ws.Range("J30").Value = 160
$ws.Range("K30").Value = 8000
$ws.Range("L30").Value = 9000
$ws.Range("M30").Value = 8500
$ws.Range("P30").Value = 425

$ws.Range("D31").Value = 44407
$ws.Range("J31").Value = 160
$ws.Range("K31").Value = 7000
$ws.Range("L31").Value = 8000
$ws.Range("M31").Value = 7500
$ws.Range("P31").Value = 375

$ws.Range("D32").Value = 44505
$ws.Range("J32").Value = 160
$ws.Range("K32").Value = 6000
$ws.Range("L32").Value = 7000
$ws.Range("M32").Value = 6500
$ws.Range("P32").Value = 325

$ws.Range("D33").Value = 44214
$ws.Range("J33").Value = 50
$ws.Range("K33").Value = 8000
$ws.Range("L33").Value = 8000
$ws.Range("M33").Value = 8000
$ws.Range("P33").Value = 400

$ws.Range("D34").Value = 44252
$ws.Range("J34").Value = 160
$ws.Range("K34").Value = 8000
$ws.Range("L34").Value = 8000
$ws.Range("M34").Value = 8000
$ws.Range("P34").Value = 400

$ws.Range("D35").Value = 44162
$ws.Range("J35").Value = 50
$ws.Range("K35").Value = 8000
$ws.Range("L35").Value = 8000
$ws.Range("M35").Value = 8000
$ws.Range("P35").Value = 400

$ws.Range("D36").Value = 44349
$ws.Range("J36").Value = 130
$ws.Range("K36").Value = 8000
$ws.Range("L36").Value = 8000
$ws.Range("M36").Value = 8000
$ws.Range("P36").Value = 400

$ws.Range("D37").Value = 44292
$ws.Range("J37").Value = 70
$ws.Range("K37").Value = 8000
$ws.Range("L37").Value = 8000
$ws.Range("M37").Value = 8000
$ws.Range("P37").Value = 400

$ws.Range("D38").Value = 44299
$ws.Range("J38").Value = 160
$ws.Range("K38").Value = 8000
$ws.Range("L38").Value = 8000
$ws.Range("M38").Value = 8000
$ws.Range("P38").Value = 400

$ws.Range("D39").Value = 44166
$ws.Range("J39").Value = 50
$ws.Range("K39").Value = 8000
$ws.Range("L39").Value = 8000
$ws.Range("M39").Value = 8000
$ws.Range("P39").Value = 400

$ws.Range("D40").Value = 44174
$ws.Range("J40").Value = 70
$ws.Range("K40").Value = 8000
$ws.Range("L40").Value = 8000
$ws.Range("M40").Value = 8000
$ws.Range("P40").Value = 400

$ws.Range("D41").Value = 44482
$ws.Range("J41").Value = 160
$ws.Range("K41").Value = 7000
$ws.Range("L41").Value = 8000
$ws.Range("M41").Value = 7500
$ws.Range("P41").Value = 375

$ws.Range("D42").Value = 44273
$ws.Range("J42").Value = 70
$ws.Range("K42").Value = 8000
$ws.Range("L42").Value = 8000
$ws.Range("M42").Value = 8000
$ws.Range("P42").Value = 400

$ws.Range("D43").Value = 44265
$ws.Range("J43").Value = 70
$ws.Range("K43").Value = 8000
$ws.Range("L43").Value = 8000
$ws.Range("M43").Value = 8000
$ws.Range("P43").Value = 400

$ws.Range("D44").Value = 44266
$ws.Range("J44").Value = 50
$ws.Range("K44").Value = 8000
$ws.Range("L44").Value = 8000
$ws.Range("M44").Value = 8000
$ws.Range("P44").Value = 400

$ws.Range("D45").Value = 44267
$ws.Range("J45").Value = 160
$ws.Range("K45").Value = 8000
$ws.Range("L45").Value = 8000
$ws.Range("M45").Value = 8000
$ws.Range("P45").Value = 400

$ws.Range("D46").Value = 44306
$ws.Range("J46").Value = 160
$ws.Range("K46").Value = 8000
$ws.Range("L46").Value = 8000
$ws.Range("M46").Value = 8000
$ws.Range("P46").Value = 400

$ws.Range("D47").Value = 44539
$ws.Range("J47").Value = 133
$ws.Range("K47").Value = 6000
$ws.Range("L47").Value = 7000
$ws.Range("M47").Value = 6504
$ws.Range("P47").Value = 325

$ws.Range("D48").Value = 44370
$ws.Range("J48").Value = 160
$ws.Range("K48").Value = 7500
$ws.Range("L48").Value = 8000
$ws.Range("M48").Value = 7750
$ws.Range("P48").Value = 388

$ws.Range("D49").Value = 44475
$ws.Range("J49").Value = 160
$ws.Range("K49").Value = 7000
$ws.Range("L49").Value = 8000
$ws.Range("M49").Value = 7500
$ws.Range("P49").Value = 375

$ws.Range("D50").Value = 44468
$ws.Range("J50").Value = 133
$ws.Range("K50").Value = 7000
$ws.Range("L50").Value = 8000
$ws.Range("M50").Value = 7504
$ws.Range("P50").Value = 375

$ws.Range("D51").Value = 44489
$ws.Range("J51").Value = 160
$ws.Range("K51").Value = 7000
$ws.Range("L51").Value = 8000
$ws.Range("M51").Value = 7500
$ws.Range("P51").Value = 375

$ws.Range("D52").Value = 44526
$ws.Range("J52").Value = 170
$ws.Range("K52").Value = 6000
$ws.Range("L52").Value = 7000
$ws.Range("M52").Value = 6500
$ws.Range("P52").Value = 325

$ws.Range("D53").Value = 44218
$ws.Range("J53").Value = 80
$ws.Range("K53").Value = 6000
$ws.Range("L53").Value = 7000
$ws.Range("M53").Value = 6625
$ws.Range("P53").Value = 331

$ws.Range("D54").Value = 44447
$ws.Range("J54").Value = 106
$ws.Range("K54").Value = 7000
$ws.Range("L54").Value = 8000
$ws.Range("M54").Value = 7500
$ws.Range("P54").Value = 375

$ws.Range("D55").Value = 44167
$ws.Range("J55").Value = 50
$ws.Range("K55").Value = 8000
$ws.Range("L55").Value = 8000
$ws.Range("M55").Value = 8000
$ws.Range("P55").Value = 400

$ws.Range("D56").Value = 44328
$ws.Range("J56").Value = 160
$ws.Range("K56").Value = 8000
$ws.Range("L56").Value = 8000
$ws.Range("M56").Value = 8000
$ws.Range("P56").Value = 400

$ws.Range("D57").Value = 44160
$ws.Range("J57").Value = 50
$ws.Range("K57").Value = 7000
$ws.Range("L57").Value = 8000
$ws.Range("M57").Value = 7600
$ws.Range("P57").Value = 380

$ws.Range("D58").Value = 44259
$ws.Range("J58").Value = 160
$ws.Range("K58").Value = 8000
$ws.Range("L58").Value = 8000
$ws.Range("M58").Value = 8000
$ws.Range("P58").Value = 400

$ws.Range("D59").Value = 44286
$ws.Range("J59").Value = 70
$ws.Range("K59").Value = 8000
$ws.Range("L59").Value = 8000
$ws.Range("M59").Value = 8000
$ws.Range("P59").Value = 400

$ws.Range("D60").Value = 44208
$ws.Range("J60").Value = 50
$ws.Range("K60").Value = 8000
$ws.Range("L60").Value = 8000
$ws.Range("M60").Value = 8000
$ws.Range("P60").Value = 400

$ws.Range("D61").Value = 44491
$ws.Range("J61").Value = 160
$ws.Range("K61").Value = 7000
$ws.Range("L61").Value = 8000
$ws.Range("M61").Value = 7500
$ws.Range("P61").Value = 375

$ws.Range("D62").Value = 44391
$ws.Range("J62").Value = 52
$ws.Range("K62").Value = 7000
$ws.Range("L62").Value = 8000
$ws.Range("M62").Value = 7500
$ws.Range("P62").Value = 375

$ws.Range("D63").Value = 44510
$ws.Range("J63").Value = 160
$ws.Range("K63").Value = 6000
$ws.Range("L63").Value = 7000
$ws.Range("M63").Value = 6500
$ws.Range("P63").Value = 325

$ws.Range("D64").Value = 44232
$ws.Range("J64").Value = 60
$ws.Range("K64").Value = 7000
$ws.Range("L64").Value = 7000
$ws.Range("M64").Value = 7000
$ws.Range("P64").Value = 350

$ws.Range("D65").Value = 44519
$ws.Range("J65").Value = 160
$ws.Range("K65").Value = 6000
$ws.Range("L65").Value = 7000
$ws.Range("M65").Value = 6500
$ws.Range("P65").Value = 325

$ws.Range("D66").Value = 44414
$ws.Range("J66").Value = 180
$ws.Range("K66").Value = 7500
$ws.Range("L66").Value = 8000
$ws.Range("M66").Value = 7750
$ws.Range("P66").Value = 388

$ws.Range("D67").Value = 44321
$ws.Range("J67").Value = 250
$ws.Range("K67").Value = 7000
$ws.Range("L67").Value = 7000
$ws.Range("M67").Value = 7000
$ws.Range("P67").Value = 350

$ws.Range("D68").Value = 44278
$ws.Range("J68").Value = 70
$ws.Range("K68").Value = 8000
$ws.Range("L68").Value = 8000
$ws.Range("M68").Value = 8000
$ws.Range("P68").Value = 400

$ws.Range("D69").Value = 44308
$ws.Range("J69").Value = 160
$ws.Range("K69").Value = 8000
$ws.Range("L69").Value = 8000
$ws.Range("M69").Value = 8000
$ws.Range("P69").Value = 400

$ws.Range("D70").Value = 44281
$ws.Range("J70").Value = 250
$ws.Range("K70").Value = 8000
$ws.Range("L70").Value = 8000
$ws.Range("M70").Value = 8000
$ws.Range("P70").Value = 400

$ws.Range("D71").Value = 44405
$ws.Range("J71").Value = 160
$ws.Range("K71").Value = 7000
$ws.Range("L71").Value = 8000
$ws.Range("M71").Value = 7500
$ws.Range("P71").Value = 375

$ws.Range("D72").Value = 44428
$ws.Range("J72").Value = 97
$ws.Range("K72").Value = 8000
$ws.Range("L72").Value = 9000
$ws.Range("M72").Value = 8505
$ws.Range("P72").Value = 425

$ws.Range("D73").Value = 44224
$ws.Range("J73").Value = 120
$ws.Range("K73").Value = 6000
$ws.Range("L73").Value = 7000
$ws.Range("M73").Value = 6667
$ws.Range("P73").Value = 333

$ws.Range("D74").Value = 44329
$ws.Range("J74").Value = 160
$ws.Range("K74").Value = 8000
$ws.Range("L74").Value = 8000
$ws.Range("M74").Value = 8000
$ws.Range("P74").Value = 400

$ws.Range("D75").Value = 44195
$ws.Range("J75").Value = 70
$ws.Range("K75").Value = 7000
$ws.Range("L75").Value = 7000
$ws.Range("M75").Value = 7000
$ws.Range("P75").Value = 350

$ws.Range("D76").Value = 44398
$ws.Range("J76").Value = 70
$ws.Range("K76").Value = 7500
$ws.Range("L76").Value = 8000
$ws.Range("M76").Value = 7750
$ws.Range("P76").Value = 388

$ws.Range("D77").Value = 44302
$ws.Range("J77").Value = 160
$ws.Range("K77").Value = 8000
$ws.Range("L77").Value = 8000
$ws.Range("M77").Value = 8000
$ws.Range("P77").Value = 400

$ws.Range("D78").Value = 44239
$ws.Range("J78").Value = 70
$ws.Range("K78").Value = 8000
$ws.Range("L78").Value = 8000
$ws.Range("M78").Value = 8000
$ws.Range("P78").Value = 400

$ws.Range("D79").Value = 44344
$ws.Range("J79").Value = 210
$ws.Range("K79").Value = 8000
$ws.Range("L79").Value = 8000
$ws.Range("M79").Value = 8000
$ws.Range("P79").Value = 400

$ws.Range("D80").Value = 44461
$ws.Range("J80").Value = 79
$ws.Range("K80").Value = 7000
$ws.Range("L80").Value = 8000
$ws.Range("M80").Value = 7494
$ws.Range("P80").Value = 375

$ws.Range("D81").Value = 44463
$ws.Range("J81").Value = 160
$ws.Range("K81").Value = 7500
$ws.Range("L81").Value = 8000
$ws.Range("M81").Value = 7750
$ws.Range("P81").Value = 388

$ws.Range("D82").Value = 44365
$ws.Range("J82").Value = 180
$ws.Range("K82").Value = 8000
$ws.Range("L82").Value = 8000
$ws.Range("M82").Value = 8000
$ws.Range("P82").Value = 400

$ws.Range("D83").Value = 44454
$ws.Range("J83").Value = 160
$ws.Range("K83").Value = 7000
$ws.Range("L83").Value = 8000
$ws.Range("M83").Value = 7500
$ws.Range("P83").Value = 375

# Append new row 84 (previously row 83 data, pushed down by the new weekly insert)
$ws.Range("A84").Value = 9
$ws.Range("B84").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C84").Value = "Metropolitana"
$ws.Range("D84").Value = 44272
$ws.Range("D84").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E84").Value = 13
$ws.Range("F84").Value = 100112005
$ws.Range("G84").Value = "Puerro"
$ws.Range("H84").Value = "Sin especificar"
$ws.Range("I84").Value = "Primera"
$ws.Range("J84").Value = 160
$ws.Range("K84").Value = 8000
$ws.Range("L84").Value = 8000
$ws.Range("M84").Value = 8000
$ws.Range("N84").Value = '$/paquete 20 unidades'
$ws.Range("O84").Value = "Provincia de Chacabuco"
$ws.Range("P84").Value = 400
$ws.Range("Q84").Value = 20
$ws.Range("R84").Value = "Hortaliza"
